$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C to fit the longer completed-task text (target stored width 60.875;
# the engine quantizes ColumnWidth to 1/7ths of a character, so 60.142857 is the
# input that lands closest, on the 60.857142857142854 bucket).
$ws.Columns.Item(3).ColumnWidth = 60.142857

# Append the new day-21 record (2019-04-16, 周二): date/time stamp, weekday,
# task description and duration, written in the same order the source file's
# shared-string table was produced in (A, B, D, C).
$ws.Range("A21").Value = "2019年4月16日13:09:33"
$ws.Range("B21").Value = "周二"
$ws.Range("D21").Value = "12:00--13:10"
$ws.Range("C21").Value = "base dao service课后练习（已整合至My_Sun工程）测试已完成"

# Leave the cursor on C21, matching the author's final selection.
$ws.Range("C21").Select()
